$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# Append a new row (row 24) to the LeetCode SQL50 tracking table.
$row = 24

$ws.Cells.Item($row, 1).Value = "619. Biggest Single Number"
$ws.Cells.Item($row, 2).Value = "Easy"
$ws.Cells.Item($row, 3).Value = "Sorting and Grouping"
$ws.Cells.Item($row, 4).Value = "Inner query. Select max(num) as num from (num group by num havingt count = 1)"

$linkText = "https://leetcode.com/problems/biggest-single-number/solutions/3839933/100-easy-fast-clean-solution/?envType=study-plan-v2&envId=top-sql-50 "
$ws.Cells.Item($row, 5).Value = $linkText
$ws.Hyperlinks.Add($ws.Cells.Item($row, 5), $linkText) | Out-Null

# Match the existing formatting used by the rest of the table:
# green "Easy" fill in column B, and the shared Hyperlink style in column E.
$ws.Cells.Item($row, 5).Style = $ws.Cells.Item($row - 1, 5).Style
$ws.Cells.Item($row, 2).Interior.Color = $ws.Cells.Item($row - 1, 2).Interior.Color

# Grow the structured table (Table2) so the new row is included.
$table = $ws.ListObjects.Item(1)
$table.Resize($ws.Range("A1:E24"))

# Update the view state (cursor position / scroll) similar to interactive editing.
$win = $ws.Application.ActiveWindow
$win.ScrollColumn = 2
$win.ScrollRow = 1
$ws.Range("E30").Select()
